# Update cryptocurrency price/volume table with latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.557.52'
$ws.Range("E2").Value = '  +4.96%  '
$ws.Range("D3").Value = '1.592.31'
$ws.Range("E3").Value = '  +1.46%  '
$ws.Range("D4").Value = "'" + '0.999'
$ws.Range("E4").Value = '  -0.71%  '
$ws.Range("D5").Value = "'" + '214.82'
$ws.Range("E5").Value = '  +1.84%  '
$ws.Range("E6").Value = '  +1.39%  '
$ws.Range("D7").Value = "'" + '0.999'
$ws.Range("E7").Value = '  -0.74%  '
$ws.Range("D8").Value = "'" + '24.05'
$ws.Range("E8").Value = '  +9.22%  '
$ws.Range("E9").Value = '  +1.28%  '
$ws.Range("E10").Value = '  +0.51%  '
$ws.Range("D11").Value = "'" + '0.0889'
$ws.Range("E11").Value = '  +2.23%  '
$ws.Range("D12").Value = '1.818.79'
$ws.Range("E12").Value = '  +1.41%  '
$ws.Range("D13").Value = '1.590.18'
$ws.Range("E13").Value = '  +1.09%  '
$ws.Range("E14").Value = '  +0.24%  '
$ws.Range("D15").Value = "'" + '0.532'
$ws.Range("E15").Value = '  +2.45%  '
$ws.Range("D16").Value = '28.519.14'
$ws.Range("E16").Value = '  +4.99%  '
$ws.Range("E17").Value = '  +2.84%  '
$ws.Range("D18").Value = "'" + '233.45'
$ws.Range("E18").Value = '  +7.81%  '
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("E20").Value = '  +1.11%  '
$ws.Range("E21").Value = '  -0.70%  '
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").Value = "'" + '9.43'
$ws.Range("E23").Value = '  +2.40%  '
$ws.Range("E24").Value = '  +0.85%  '
$ws.Range("D25").Value = "'" + '151.81'
$ws.Range("E25").Value = '  -1.40%  '
$ws.Range("D26").Value = "'" + '15.35'
$ws.Range("E26").Value = '  +1.83%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("E28").Value = '  +0.97%  '
$ws.Range("E29").Value = '  -0.72%  '
$ws.Range("E30").Value = '  +0.51%  '
$ws.Range("E31").Value = '  +0.37%  '
$ws.Range("E32").Value = '  +0.36%  '
$ws.Range("D33").Value = "'" + '3.14'
$ws.Range("E33").Value = '  -0.45%  '
$ws.Range("D34").Value = '1.419.20'
$ws.Range("E34").Value = '  -2.15%  '
$ws.Range("E35").Value = '  -0.76%  '
$ws.Range("D36").Value = "'" + '1.06'
$ws.Range("E36").Value = '  -5.83%  '
$ws.Range("E37").Value = '  -0.49%  '
$ws.Range("E38").Value = '  +0.32%  '
$ws.Range("E39").Value = '  +9.28%  '
$ws.Range("D40").Value = "'" + '0.543'
$ws.Range("E40").Value = '  +2.09%  '
$ws.Range("D41").Value = "'" + '0.815'
$ws.Range("E41").Value = '  +0.68%  '
$ws.Range("D42").Value = "'" + '5.75'
$ws.Range("E42").Value = '  -1.74%  '
$ws.Range("D43").Value = "'" + '0.999'
$ws.Range("E43").Value = '  -0.74%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").Value = "'" + '0.980'
$ws.Range("E44").Value = '  -2.43%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").Value = "'" + '1.83'
$ws.Range("E45").Value = '  +6.44%  '
$ws.Range("D46").Value = "'" + '64.77'
$ws.Range("E46").Value = '  +0.66%  '
$ws.Range("D47").Value = '1.730.67'
$ws.Range("E47").Value = '  +1.57%  '
$ws.Range("D48").Value = "'" + '87.97'
$ws.Range("E48").Value = '  +2.24%  '
$ws.Range("D49").Value = '0.0₆0103'
$ws.Range("E49").Value = '  -0.25%  '
$ws.Range("E50").Value = '  -0.18%  '
$ws.Range("D51").Value = "'" + '40.02'
$ws.Range("E51").Value = '  +17.95%  '
